# Daily attendance processing - 2025-12-24 09:35:55
#
# Normalizes the "Recorded By" (column G) entries so that automated /
# service accounts are listed after the human-facing address rather than
# before it, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".
# Only applies to the two accounts affected by the recent recording-agent
# change (backup@backdoor.com / dnasr281@gmail.com); other combinations
# (e.g. admin@admin.com, or entries where System is not first) are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) { continue }

    $parts = $val.Split(",")
    if ($parts.Count -lt 2) { continue }
    if ($parts[0].Trim() -ne "System") { continue }

    $target = $parts[1].Trim()
    if ($target -ne "backup@backdoor.com" -and $target -ne "dnasr281@gmail.com") {
        continue
    }

    $rest = ""
    for ($i = 2; $i -lt $parts.Count; $i++) {
        $rest = $rest + ", " + $parts[$i].Trim()
    }

    $newVal = $target + ", System" + $rest
    $cell.Value = $newVal
}
